$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===== MODIFIED CELLS =====
$ws.Cells.Item(485, 6).Value = 14219  # was 14213
$ws.Cells.Item(492, 6).Value = 14432  # was 14425
$ws.Cells.Item(512, 6).Value = 8688  # was 8683
$ws.Cells.Item(523, 6).Value = 10377  # was 10369
$ws.Cells.Item(552, 6).Value = 15626  # was 15621
$ws.Cells.Item(567, 6).Value = 23494  # was 23492
$ws.Cells.Item(569, 6).Value = 32402  # was 32387
$ws.Cells.Item(572, 6).Value = 33341  # was 33339
$ws.Cells.Item(573, 6).Value = 26992  # was 26989
$ws.Cells.Item(574, 6).Value = 23355  # was 23352
$ws.Cells.Item(575, 6).Value = 26072  # was 25976
$ws.Cells.Item(575, 7).Value = 394  # was 386
$ws.Cells.Item(579, 6).Value = 32609  # was 32604
$ws.Cells.Item(580, 6).Value = 28838  # was 28835
$ws.Cells.Item(581, 6).Value = 27033  # was 27026
$ws.Cells.Item(582, 6).Value = 25852  # was 25849
$ws.Cells.Item(583, 6).Value = 29251  # was 29241
$ws.Cells.Item(586, 6).Value = 33591  # was 33587
$ws.Cells.Item(587, 6).Value = 28183  # was 28177
$ws.Cells.Item(588, 6).Value = 25357  # was 25350
$ws.Cells.Item(589, 6).Value = 25395  # was 25378
$ws.Cells.Item(590, 6).Value = 29231  # was 29224
$ws.Cells.Item(593, 6).Value = 36919  # was 36913
$ws.Cells.Item(593, 7).Value = 1187  # was 1186
$ws.Cells.Item(594, 6).Value = 29702  # was 29699
$ws.Cells.Item(595, 6).Value = 27237  # was 27216
$ws.Cells.Item(596, 6).Value = 29040  # was 28730
$ws.Cells.Item(596, 7).Value = 941  # was 932
$ws.Cells.Item(597, 6).Value = 29816  # was 29806
$ws.Cells.Item(598, 6).Value = 15020  # was 15013
$ws.Cells.Item(599, 6).Value = 16465  # was 16321
$ws.Cells.Item(599, 7).Value = 870  # was 860
$ws.Cells.Item(600, 6).Value = 39608  # was 39293
$ws.Cells.Item(600, 7).Value = 1660  # was 1650
$ws.Cells.Item(601, 6).Value = 31315  # was 30768
$ws.Cells.Item(601, 7).Value = 1320  # was 1301
$ws.Cells.Item(602, 6).Value = 29383  # was 22199
$ws.Cells.Item(602, 7).Value = 1271  # was 961

# ===== NEW ROWS =====
# Row 603 (new)
$ws.Cells.Item(603, 1).Value = 44497
$ws.Cells.Item(603, 2).Value = 474595
$ws.Cells.Item(603, 3).Value = 17338
$ws.Cells.Item(603, 4).Value = 4587
$ws.Cells.Item(603, 5).Value = 13000
$ws.Cells.Item(603, 6).Value = 31258
$ws.Cells.Item(603, 7).Value = 1515

# Row 604 (new)
$ws.Cells.Item(604, 1).Value = 44498
$ws.Cells.Item(604, 2).Value = 479737
$ws.Cells.Item(604, 3).Value = 17638
$ws.Cells.Item(604, 4).Value = 5142
$ws.Cells.Item(604, 5).Value = 13018
$ws.Cells.Item(604, 6).Value = 26392
$ws.Cells.Item(604, 7).Value = 1410

# Row 605 (new)
$ws.Cells.Item(605, 1).Value = 44499
$ws.Cells.Item(605, 2).Value = 483773
$ws.Cells.Item(605, 3).Value = 14382
$ws.Cells.Item(605, 4).Value = 4036
$ws.Cells.Item(605, 5).Value = 13034
$ws.Cells.Item(605, 6).Value = 13022
$ws.Cells.Item(605, 7).Value = 914

# Row 606 (new)
$ws.Cells.Item(606, 1).Value = 44500
$ws.Cells.Item(606, 2).Value = 485629
$ws.Cells.Item(606, 3).Value = 6733
$ws.Cells.Item(606, 4).Value = 1856
$ws.Cells.Item(606, 5).Value = 13045
$ws.Cells.Item(606, 6).Value = 11788
$ws.Cells.Item(606, 7).Value = 1071

# Row 607 (new)
$ws.Cells.Item(607, 1).Value = 44501
$ws.Cells.Item(607, 2).Value = 487916
$ws.Cells.Item(607, 3).Value = 9682
$ws.Cells.Item(607, 4).Value = 2287
$ws.Cells.Item(607, 5).Value = 13076
$ws.Cells.Item(607, 6).Value = 7418
$ws.Cells.Item(607, 7).Value = 676

